$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "57.414.77"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "2.362.06"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.43%  "
Set-TextValue "D5" "520.69"
$ws.Range("E5").Value = "  +0.37%  "
Set-TextValue "D6" "135.90"
$ws.Range("E6").Value = "  +1.13%  "
Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  -0.04%  "
Set-TextValue "D8" "0.540"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  -0.96%  "
Set-TextValue "D10" "5.46"
$ws.Range("E10").Value = "  +5.47%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  -0.41%  "
Set-TextValue "D13" "24.41"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "2.787.31"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "57.474.29"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "2.369.92"
$ws.Range("E17").Value = "  +1.88%  "
Set-TextValue "D18" "10.62"
$ws.Range("E18").Value = "  +0.71%  "
Set-TextValue "D19" "330.40"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("E20").Value = "  -0.86%  "
Set-TextValue "D21" "6.72"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("E22").Value = "  -0.07%  "
Set-TextValue "D23" "61.31"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D24" "0.166"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D25" "8.68"
$ws.Range("E25").Value = "  +12.69%  "
Set-TextValue "D26" "0.995"
$ws.Range("E26").Value = "  -0.17%  "
Set-TextValue "D27" "1.34"
$ws.Range("E27").Value = "  +11.14%  "
$ws.Range("D28").Value = "0.0₃0746"
$ws.Range("E28").Value = "  +0.88%  "
Set-TextValue "D29" "168.37"
$ws.Range("E29").Value = "  -2.32%  "
Set-TextValue "D30" "1.70"
$ws.Range("E30").Value = "  +1.25%  "
Set-TextValue "D31" "6.31"
$ws.Range("E31").Value = "  +0.56%  "
Set-TextValue "D32" "18.59"
$ws.Range("E32").Value = "  +1.30%  "
Set-TextValue "D34" "1.31"
$ws.Range("E34").Value = "  +3.69%  "
Set-TextValue "D35" "0.995"
$ws.Range("E35").Value = "  +0.33%  "
Set-TextValue "D36" "0.923"
$ws.Range("E36").Value = "  -2.49%  "
Set-TextValue "D37" "4.05"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  +6.36%  "
Set-TextValue "D39" "38.69"
$ws.Range("E39").Value = "  +2.96%  "
Set-TextValue "D40" "150.25"
$ws.Range("E40").Value = "  +6.83%  "
$ws.Range("E41").Value = "  +1.39%  "
Set-TextValue "D42" "3.67"
$ws.Range("E42").Value = "  +1.68%  "
Set-TextValue "D43" "5.36"
$ws.Range("E43").Value = "  +3.94%  "
Set-TextValue "D44" "283.93"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("E45").Value = "  +1.27%  "
Set-TextValue "D46" "0.0510"
$ws.Range("E46").Value = "  +0.00%  "
Set-TextValue "D47" "0.565"
$ws.Range("E47").Value = "  +0.31%  "
Set-TextValue "D48" "18.30"
$ws.Range("E48").Value = "  +6.16%  "
Set-TextValue "D49" "0.0220"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D50" "0.388"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "17.73"
$ws.Range("E51").Value = "  +4.58%  "
